# Apply the "Add files via upload" edit:
#  - Add two new daily-message rows (8 and 9), mirroring the existing
#    daily-message row pairs (A/D columns repeat the same two shared
#    strings; B holds a new 5-minute time window; C stays blank).
#  - Update the hourly-message row (25) to a new time window.
#  - Update the sheet view (scrolled position + active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 8: mirrors rows 2/4/6 (A = prayer text, D = "ежедневный посыл") ---
$ws.Range("A8").Value = $ws.Range("A6").Value2
$ws.Range("B8").Value = "18:20 - 18:24"
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = $ws.Range("D6").Value2

# --- New row 9: mirrors rows 3/5/7 ---
$ws.Range("A9").Value = $ws.Range("A7").Value2
$ws.Range("B9").Value = "18:25 - 18:29"
$ws.Range("D9").Value = $ws.Range("D7").Value2

# --- Row 25 (hourly message): new time window replaces the old one ---
$ws.Range("B25").Value = "18:35 - 18:39"

# --- Row heights for the new rows, matching rows 6/7 exactly ---
$ws.Rows.Item(8).RowHeight = $ws.Rows.Item(6).RowHeight
$ws.Rows.Item(9).RowHeight = $ws.Rows.Item(7).RowHeight

# --- Style the new cells like their mirrored counterparts ---
$ws.Range("A8").Style = $ws.Range("A6").Style
$ws.Range("C8").Style = $ws.Range("C6").Style

# --- Update sheet view: scrolled to row 11, selection on C30 ---
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("C30").Select()
